$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New row of data appended at row 95, matching the existing table structure
$rowValues = @(44334, 515, 3001, 900, 6077, 253, 2079, 1000, 28500, 0, 0, 0, 3674, 184, 1099, 43330)

for ($i = 0; $i -lt $rowValues.Length; $i++) {
    $col = $i + 1
    $ws.Cells.Item(95, $col).Value = $rowValues[$i]
}

# Column A uses the same date number-format style as the rest of column A (yyyy-mm-dd)
$ws.Range("A95").NumberFormat = "yyyy-mm-dd"
